# Update column F (dSF) values on the active worksheet per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 4
    5  = -2
    6  = 7
    7  = -6
    8  = -1
    9  = 4
    12 = -4
    14 = -4
    16 = 4
    17 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
